# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment schedule"
# sheet, shifting the old N/O/P ("Late" / spare-heading / "Outstanding")
# columns one place to the right, then make that sheet the active tab
# (picking up the selection that was left on P7) instead of "Transactions".

$wb = $excel.ActiveWorkbook

$wsRepayment = $wb.Worksheets.Item("Repayment schedule")

# Make "Repayment schedule" the active sheet (was "Transactions"); this
# also clears the tabSelected flag that was on "Transactions".
$wsRepayment.Activate()

# Insert a new column at N; everything from N onward shifts right by one.
# The new blank column inherits the width of the column to its left (M).
$mColumnWidth = $wsRepayment.Range("M1").ColumnWidth
$wsRepayment.Columns("N").Insert()
$wsRepayment.Range("N1").ColumnWidth = $mColumnWidth

# Leave the same selection Excel left after the insert.
$wsRepayment.Range("P7").Select()
